$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add hyphen separators to the three "cover page" label strings ---
$ws.Range("A5").Value = "C* - A* - B:* CAB"
$ws.Range("A6").Value = "B* - A* - D:* BAD"
$ws.Range("A7").Value = "D* - A* - B:* DAB"

# --- B5 ("cab.jpg") now picks up the same non-wrapping text style as the
#     rest of the label column (A2/A4/A5/A6/A7) instead of the header style ---
$ws.Range("B5").Font.Name = "Arial"
$ws.Range("B5").Font.Size = 10
$ws.Range("B5").Font.ThemeColor = 1
$ws.Range("B5").WrapText = $false

# --- Cover page view tweaks: zoomed in, with the selection resting on A7 ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 209
$ws.Range("A7").Select()
